# Auto-generated script applying the crypto price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number must be forced to a
# Text number format first, otherwise Excel would silently convert the
# cell from a text/inline-string cell into a numeric cell.
$textCells = @("D5","D6","D13","D18","D19","D21","D22","D24","D25","D27","D28","D31","D34","D38","D42","D43","D44","D46","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = '69.148.98'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '3.749.35'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '601.66'
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").Value = '167.04'
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("D7").Value = '3.746.78'
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("E10").Value = '  +3.67%  '
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Value = '37.92'
$ws.Range("E14").Value = '  +1.61%  '
$ws.Range("D15").Value = '4.375.29'
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").Value = '3.765.02'
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("D17").Value = '69.151.65'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").Value = '7.34'
$ws.Range("E18").Value = '  +1.52%  '
$ws.Range("D19").Value = '17.37'
$ws.Range("E19").Value = '  +0.97%  '
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("D21").Value = '11.05'
$ws.Range("E21").Value = '  +8.00%  '
$ws.Range("D22").Value = '493.05'
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("E23").Value = '  +0.55%  '
$ws.Range("D24").Value = '0.0000152'
$ws.Range("E24").Value = '  +8.25%  '
$ws.Range("D25").Value = '84.87'
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("D27").Value = '12.28'
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("D28").Value = '10.07'
$ws.Range("E28").Value = '  -0.72%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +0.94%  '
$ws.Range("D31").Value = '8.14'
$ws.Range("E31").Value = '  +1.92%  '
$ws.Range("E32").Value = '  +1.57%  '
$ws.Range("D33").Value = '3.895.41'
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("D34").Value = '31.45'
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("D35").Value = '3.682.70'
$ws.Range("E35").Value = '  +0.47%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = '1.02'
$ws.Range("E38").Value = '  +0.43%  '
$ws.Range("E39").Value = '  +2.54%  '
$ws.Range("E40").Value = '  +3.31%  '
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("D42").Value = '3.04'
$ws.Range("E42").Value = '  +6.32%  '
$ws.Range("D43").Value = '48.79'
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("D44").Value = '425.19'
$ws.Range("E44").Value = '  -2.52%  '
$ws.Range("E45").Value = '  -0.47%  '
$ws.Range("D46").Value = '8.45'
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("E48").Value = '  -1.07%  '
$ws.Range("D49").Value = '141.68'
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").Value = '2.789.71'
$ws.Range("E50").Value = '  +1.66%  '
$ws.Range("E51").Value = '  +0.19%  '
